$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.484390497207642
$ws.Range("B1").Value = 3.599098205566406
$ws.Range("C1").Value = 2.916921138763428
$ws.Range("D1").Value = 1.370304226875305
$ws.Range("E1").Value = 0.7716301083564758
